# Staging.MilestoneType.xlsx edit
#
# Per the commit's xml diff, the two header-row shared strings used by
# A2/B2 ("MilestoneTypeID" and "Code") swap places, while the cell
# references (A2 -> shared string #1, B2 -> shared string #2) stay put.
# The net, observable effect is simply that the header row's first two
# cells swap their displayed text:
#   A2: "MilestoneTypeID" -> "Code"
#   B2: "Code"            -> "MilestoneTypeID"
# C2 ("Name") is untouched.
#
# (The diff's other hunks - the workbookView window-size sentinel, the
# worksheet's VBA codeName, and dropping the B/C <col> width overrides -
# are cosmetic/IDE-state artifacts with no corresponding writable
# property on the Excel object model exposed here, so they are left
# alone rather than risk corrupting the sheet.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Code"
$ws.Range("B2").Value = "MilestoneTypeID"
